$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AMSIN")
$ws1.Rows("14:14").Insert()
$ws1.Cells.Item(14, 1).Formula = "=""2024-07-15"""
$ws1.Range("A14").Copy()
$ws1.Range("A14").PasteSpecial(-4163)
$ws1.Cells.Item(14, 2).Value = 45488.79296148148
$ws1.Cells.Item(14, 3).Formula = "=""193livessd"""
$ws1.Range("C14").Copy()
$ws1.Range("C14").PasteSpecial(-4163)
$ws1.Cells.Item(14, 4).Value = 45
$ws1.Cells.Item(14, 5).Value = 56
$ws1.Cells.Item(14, 6).Value = -11
$ws1.Cells.Item(14, 7).Value = 5.79
Write-Output "row 14 on $ws1 filled ok"
$ws1.Rows("15:15").Delete()
$ws1.Rows("15:15").Insert()
$ws1.Cells.Item(15, 1).Formula = "=""2024-07-16"""
$ws1.Range("A15").Copy()
$ws1.Range("A15").PasteSpecial(-4163)
$ws1.Cells.Item(15, 2).Value = 45489.46580332176
$ws1.Cells.Item(15, 3).Formula = "=""193fghj"""
$ws1.Range("C15").Copy()
$ws1.Range("C15").PasteSpecial(-4163)
$ws1.Cells.Item(15, 4).Value = 60
$ws1.Cells.Item(15, 5).Value = 0
$ws1.Cells.Item(15, 6).Value = 60
$ws1.Cells.Item(15, 7).Value = 0.09
Write-Output "row 15 on $ws1 filled ok"
$ws1.Rows("16:16").Insert()
$ws1.Cells.Item(16, 1).Formula = "=""2024-07-16"""
$ws1.Range("A16").Copy()
$ws1.Range("A16").PasteSpecial(-4163)
$ws1.Cells.Item(16, 2).Value = 45489.46703162037
$ws1.Cells.Item(16, 3).Formula = "=""193dfghj"""
$ws1.Range("C16").Copy()
$ws1.Range("C16").PasteSpecial(-4163)
$ws1.Cells.Item(16, 4).Value = 60
$ws1.Cells.Item(16, 5).Value = 0
$ws1.Cells.Item(16, 6).Value = 60
$ws1.Cells.Item(16, 7).Value = 0.07
Write-Output "row 16 on $ws1 filled ok"
$ws1.Rows("17:17").Insert()
$ws1.Cells.Item(17, 1).Formula = "=""2024-07-16"""
$ws1.Range("A17").Copy()
$ws1.Range("A17").PasteSpecial(-4163)
$ws1.Cells.Item(17, 2).Value = 45489.47452182871
$ws1.Cells.Item(17, 3).Formula = "=""aas22"""
$ws1.Range("C17").Copy()
$ws1.Range("C17").PasteSpecial(-4163)
$ws1.Cells.Item(17, 4).Value = 60
$ws1.Cells.Item(17, 5).Value = 20
$ws1.Cells.Item(17, 6).Value = 40
$ws1.Cells.Item(17, 7).Value = 0.68
Write-Output "row 17 on $ws1 filled ok"
$ws1.Rows("18:18").Insert()
$ws1.Cells.Item(18, 1).Formula = "=""2024-07-16"""
$ws1.Range("A18").Copy()
$ws1.Range("A18").PasteSpecial(-4163)
$ws1.Cells.Item(18, 2).Value = 45489.47601700231
$ws1.Cells.Item(18, 3).Formula = "=""193jjer"""
$ws1.Range("C18").Copy()
$ws1.Range("C18").PasteSpecial(-4163)
$ws1.Cells.Item(18, 4).Value = 60
$ws1.Cells.Item(18, 5).Value = 83
$ws1.Cells.Item(18, 6).Value = -23
$ws1.Cells.Item(18, 7).Value = 8.19
Write-Output "row 18 on $ws1 filled ok"
$ws1.Rows("19:19").Insert()
$ws1.Cells.Item(19, 1).Formula = "=""2024-07-16"""
$ws1.Range("A19").Copy()
$ws1.Range("A19").PasteSpecial(-4163)
$ws1.Cells.Item(19, 2).Value = 45489.54144825231
$ws1.Cells.Item(19, 3).Formula = "=""193vinodds"""
$ws1.Range("C19").Copy()
$ws1.Range("C19").PasteSpecial(-4163)
$ws1.Cells.Item(19, 4).Value = 60
$ws1.Cells.Item(19, 5).Value = 113
$ws1.Cells.Item(19, 6).Value = -53
$ws1.Cells.Item(19, 7).Value = 5.05
Write-Output "row 19 on $ws1 filled ok"
$ws1.Rows("20:20").Insert()
$ws1.Cells.Item(20, 1).Formula = "=""2024-07-16"""
$ws1.Range("A20").Copy()
$ws1.Range("A20").PasteSpecial(-4163)
$ws1.Cells.Item(20, 2).Value = 45489.55096924768
$ws1.Cells.Item(20, 3).Formula = "=""193ert"""
$ws1.Range("C20").Copy()
$ws1.Range("C20").PasteSpecial(-4163)
$ws1.Cells.Item(20, 4).Value = 116
$ws1.Cells.Item(20, 5).Value = 33
$ws1.Cells.Item(20, 6).Value = 83
$ws1.Cells.Item(20, 7).Value = 1.75
Write-Output "row 20 on $ws1 filled ok"
$ws1.Rows("21:21").Insert()
$ws1.Cells.Item(21, 1).Formula = "=""2024-07-16"""
$ws1.Range("A21").Copy()
$ws1.Range("A21").PasteSpecial(-4163)
$ws1.Cells.Item(21, 2).Value = 45489.5618215162
$ws1.Cells.Item(21, 3).Formula = "=""193sdsdd"""
$ws1.Range("C21").Copy()
$ws1.Range("C21").PasteSpecial(-4163)
$ws1.Cells.Item(21, 4).Value = 116
$ws1.Cells.Item(21, 5).Value = 34
$ws1.Cells.Item(21, 6).Value = 82
$ws1.Cells.Item(21, 7).Value = 1.46
Write-Output "row 21 on $ws1 filled ok"
$ws1.Rows("22:22").Insert()
$ws1.Cells.Item(22, 1).Formula = "=""2024-07-16"""
$ws1.Range("A22").Copy()
$ws1.Range("A22").PasteSpecial(-4163)
$ws1.Cells.Item(22, 2).Value = 45489.57467228009
$ws1.Cells.Item(22, 3).Value = ""
$ws1.Cells.Item(22, 4).Value = 116
$ws1.Cells.Item(22, 5).Value = 34
$ws1.Cells.Item(22, 6).Value = 82
$ws1.Cells.Item(22, 7).Value = 1.23
Write-Output "row 22 on $ws1 filled ok"
$ws1.Rows("23:23").Insert()
$ws1.Cells.Item(23, 1).Formula = "=""2024-07-16"""
$ws1.Range("A23").Copy()
$ws1.Range("A23").PasteSpecial(-4163)
$ws1.Cells.Item(23, 2).Value = 45489.57619274306
$ws1.Cells.Item(23, 3).Formula = "=""193lattest"""
$ws1.Range("C23").Copy()
$ws1.Range("C23").PasteSpecial(-4163)
$ws1.Cells.Item(23, 4).Value = 116
$ws1.Cells.Item(23, 5).Value = 114
$ws1.Cells.Item(23, 6).Value = 2
$ws1.Cells.Item(23, 7).Value = 4.46
Write-Output "row 23 on $ws1 filled ok"
$ws1.Rows("24:24").Insert()
$ws1.Cells.Item(24, 1).Formula = "=""2024-07-16"""
$ws1.Range("A24").Copy()
$ws1.Range("A24").PasteSpecial(-4163)
$ws1.Cells.Item(24, 2).Value = 45489.58621821759
$ws1.Cells.Item(24, 3).Formula = "=""193vinoddd"""
$ws1.Range("C24").Copy()
$ws1.Range("C24").PasteSpecial(-4163)
$ws1.Cells.Item(24, 4).Value = 115
$ws1.Cells.Item(24, 5).Value = 112
$ws1.Cells.Item(24, 6).Value = 3
$ws1.Cells.Item(24, 7).Value = 4.96
Write-Output "row 24 on $ws1 filled ok"

$ws2 = $wb.Worksheets.Item("AMS")
$ws2.Rows("11:11").Insert()
$ws2.Cells.Item(11, 1).Formula = "=""2024-07-16"""
$ws2.Range("A11").Copy()
$ws2.Range("A11").PasteSpecial(-4163)
$ws2.Cells.Item(11, 2).Value = 45489.59596949074
$ws2.Cells.Item(11, 3).Formula = "=""193livee"""
$ws2.Range("C11").Copy()
$ws2.Range("C11").PasteSpecial(-4163)
$ws2.Cells.Item(11, 4).Value = 115
$ws2.Cells.Item(11, 5).Value = 114
$ws2.Cells.Item(11, 6).Value = 1
$ws2.Cells.Item(11, 7).Value = 5.38
Write-Output "row 11 on $ws2 filled ok"
$ws2.Rows("12:12").Insert()
$ws2.Range("A12").ClearFormats()
$ws2.Range("C12:G12").ClearFormats()
$ws2.Cells.Item(12, 1).Formula = "=""2024-07-16"""
$ws2.Range("A12").Copy()
$ws2.Range("A12").PasteSpecial(-4163)
$ws2.Cells.Item(12, 2).Value = 45489.60234038021
$ws2.Cells.Item(12, 3).Formula = "=""194ddsds"""
$ws2.Range("C12").Copy()
$ws2.Range("C12").PasteSpecial(-4163)
$ws2.Cells.Item(12, 4).Value = 115
$ws2.Cells.Item(12, 5).Value = 114
$ws2.Cells.Item(12, 6).Value = 1
$ws2.Cells.Item(12, 7).Value = 4.11
Write-Output "row 12 on $ws2 filled ok"
